# In Progress - WorkOrder Loading
# Clear the stale placeholder/group-row values out of row 2 on the WOCENTER
# sheet (these were internal "[GROUPROW]" / "188:x" / "189:x" / "8065:x"
# marker strings left over from a prior data pull) while keeping the row's
# formatting intact, and move the active selection up to A2 to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WOCENTER")

$ws.Range("A2:AI2").Value = $null

$null = $ws.Range("A2").Select()
